$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.095493383854034164
$ws.Range("B1").Value = 0.095441033559694688
$ws.Range("A2").Value = -0.073336552294362534
$ws.Range("B2").Value = 0.073205318536783359
$ws.Range("A3").Value = -0.023502122991105878
$ws.Range("B3").Value = 0.023469383855728054
$ws.Range("A4").Value = -0.015469383902997791
$ws.Range("B4").Value = 0.014790427645623438
$ws.Range("A5").Value = -0.011790427666843684
$ws.Range("B5").Value = 0.0094424824106278393
$ws.Range("A6").Value = -0.052671755124517716
$ws.Range("B6").Value = 0.052078346472107739
$ws.Range("A7").Value = -0.042078346533521938
$ws.Range("B7").Value = 0.041928763784834455
$ws.Range("A8").Value = -0.031928763848966035
$ws.Range("B8").Value = 0.031647543714333093
$ws.Range("A9").Value = -0.029647543739519833
$ws.Range("B9").Value = 0.029411294339610983
$ws.Range("A10").Value = -0.027411294368564043
$ws.Range("B10").Value = 0.027395593536082785
$ws.Range("A11").Value = -0.024395593570852192
$ws.Range("B11").Value = 0.024369021303018989
$ws.Range("A12").Value = -0.020869021341222371
$ws.Range("B12").Value = 0.020672596189982961
$ws.Range("A13").Value = -0.017172596231951331
$ws.Range("B13").Value = 0.01708316511825192
$ws.Range("A14").Value = -0.0090831651854044182
$ws.Range("B14").Value = 0.0090539707435768335
$ws.Range("A15").Value = -0.008053970773790553
$ws.Range("B15").Value = 0.0080351190271894524
$ws.Range("A16").Value = -0.0060351190633851637
$ws.Range("B16").Value = 0.0060035515422800501
$ws.Range("A17").Value = -0.0040035515794105692
$ws.Range("B17").Value = 0.0039999999519269025
$ws.Range("A18").Value = -0.01610448129811104
$ws.Range("B18").Value = 0.016091341531495829
$ws.Range("A19").Value = -0.012091341553689627
$ws.Range("B19").Value = 0.012016650151823072
$ws.Range("A20").Value = -0.0080166501758824893
$ws.Range("B20").Value = 0.0080056598992150896
$ws.Range("A21").Value = -0.004005659923537408
$ws.Range("B21").Value = 0.0039999999754920523
$ws.Range("A22").Value = -0.045703195570753863
$ws.Range("B22").Value = 0.04549254939704106
$ws.Range("A23").Value = -0.04049254942931757
$ws.Range("B23").Value = 0.040097727355276724
$ws.Range("A24").Value = -0.02009772747167915
$ws.Range("B24").Value = 0.019999999881900266
$ws.Range("A25").Value = -0.011034051450462101
$ws.Range("B25").Value = 0.010994386747578133
$ws.Range("A26").Value = -0.0084943867677207407
$ws.Range("B26").Value = 0.0084464523800100721
$ws.Range("A27").Value = -0.005946452400388047
$ws.Range("B27").Value = 0.0056813821155974331
$ws.Range("A28").Value = -0.0036813821342454034
$ws.Range("B28").Value = 0.0035181843837710502
$ws.Range("A29").Value = 0.0034818155703426612
$ws.Range("B29").Value = -0.0035195686307538665
$ws.Range("A30").Value = 0.0035001432508638963
$ws.Range("B30").Value = -0.0035928332166994537
$ws.Range("A31").Value = -0.014022550177035242
$ws.Range("B31").Value = 0.014001239799087628
$ws.Range("A32").Value = -0.0040012398610294042
$ws.Range("B32").Value = 0.0039999999701443301

$ws.Columns.Item(1).ColumnWidth = 15.42578125
$ws.Columns.Item(2).ColumnWidth = 15.42578125
